$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shrink row 19 height (was 45, now 30) to match new content
$ws.Rows.Item(19).RowHeight = 30

# Row 35: Jan 31 15:00 to 16:00
$ws.Range("A35").Value = "Jan 31 15:00 to 16:00"
$ws.Range("B35").Value = "Build sample multiple linear regression model."
$ws.Range("C35").Value = "Infimetrics"

# Row 36: Jan 31 16:00 to 17:00
$ws.Range("A36").Value = "Jan 31 16:00 to 17:00"
$ws.Range("B36").Value = "Build linear regression model using sklearn.linear regresion"
$ws.Range("C36").Value = "Infimetrics"

# Row 37: Jan 31 17:00 to 18:00
$ws.Range("A37").Value = "Jan 31 17:00 to 18:00"
$ws.Range("B37").Value = "Seached about other predictive model. Undertanding features of data."
$ws.Range("C37").Value = "Infimetrics"

# Row 38: Jan 31 18:00 to 19:00
$ws.Range("A38").Value = "Jan 31 18:00 to 19:00"
$ws.Range("B38").Value = "Applied MLP neural network. Results were all same values. Need to fix this problem because further neural network strategy will be used for prediction. "
$ws.Range("C38").Value = "Infimetrics"

# Apply style formatting consistent with existing rows:
#  A/C columns: center/center, no wrap (style 1)
#  B column rows 35-37 (short single-line text): left/center, no wrap (style 3)
#  B column row 38 (long text): left/center, wrapText (style 2)
$ws.Range("A35:A37").HorizontalAlignment = -4108
$ws.Range("A35:A37").VerticalAlignment = -4108
$ws.Range("A35:A37").WrapText = $false

$ws.Range("B35:B37").HorizontalAlignment = -4131
$ws.Range("B35:B37").VerticalAlignment = -4108
$ws.Range("B35:B37").WrapText = $false

$ws.Range("C35:C37").HorizontalAlignment = -4108
$ws.Range("C35:C37").VerticalAlignment = -4108
$ws.Range("C35:C37").WrapText = $false

$ws.Range("A38").HorizontalAlignment = -4108
$ws.Range("A38").VerticalAlignment = -4108
$ws.Range("A38").WrapText = $false

$ws.Range("B38").HorizontalAlignment = -4131
$ws.Range("B38").VerticalAlignment = -4108
$ws.Range("B38").WrapText = $true

$ws.Range("C38").HorizontalAlignment = -4108
$ws.Range("C38").VerticalAlignment = -4108
$ws.Range("C38").WrapText = $false

$ws.Rows.Item(38).RowHeight = 45

# Update selection / view to match final state
$ws.Range("D38").Select()
